$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 261-262, pushing the existing rows
# (261-265) down to (263-267), preserving their content/formatting.
$ws.Range("A261:A262").EntireRow.Insert()

# Row 261 - new weekly record (week of 44448)
$ws.Range("A261").Value = 3
$ws.Range("B261").Value = "Femacal de La Calera"
$ws.Range("C261").Value = "Coquimbo"
$ws.Range("D261").Value = 44448
$ws.Range("E261").Value = 5
$ws.Range("F261").Value = 100112037
$ws.Range("G261").Value = "Cebollín"
$ws.Range("H261").Value = "Sin especificar"
$ws.Range("I261").Value = "Primera"
$ws.Range("J261").Value = 280
$ws.Range("K261").Value = 3500
$ws.Range("L261").Value = 3800
$ws.Range("M261").Value = 3639
$ws.Range("N261").Value = "$/paquete 36 unidades"
$ws.Range("O261").Value = "Provincia de Quillota"
$ws.Range("P261").Value = 101
$ws.Range("Q261").Value = 36
$ws.Range("R261").Value = "Hortaliza"

# Row 262 - new weekly record (week of 44448)
$ws.Range("A262").Value = 3
$ws.Range("B262").Value = "Femacal de La Calera"
$ws.Range("C262").Value = "Coquimbo"
$ws.Range("D262").Value = 44448
$ws.Range("E262").Value = 5
$ws.Range("F262").Value = 100112037
$ws.Range("G262").Value = "Cebollín"
$ws.Range("H262").Value = "Sin especificar"
$ws.Range("I262").Value = "Segunda"
$ws.Range("J262").Value = 150
$ws.Range("K262").Value = 2500
$ws.Range("L262").Value = 2500
$ws.Range("M262").Value = 2500
$ws.Range("N262").Value = "$/paquete 36 unidades"
$ws.Range("O262").Value = "Provincia de Quillota"
$ws.Range("P262").Value = 69
$ws.Range("Q262").Value = 36
$ws.Range("R262").Value = "Hortaliza"
